# hodiny.xlsx update — add new time-tracking entries for 2.8., 3.8. and 12.8.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 126 — 2.8.
$ws.Range("B126").Value = "2.8."
$ws.Range("C126").Value = 0.50347222222222221
$ws.Range("D126").Value = 0.53749999999999998

# Row 127 — 3.8.
$ws.Range("B127").Value = "3.8."
$ws.Range("C127").Value = 0.29305555555555557
$ws.Range("D127").Value = 0.3263888888888889

# Row 128 — 12.8. (first entry of the day)
$ws.Range("B128").Value = "12.8."
$ws.Range("C128").Value = 0.53541666666666665
$ws.Range("D128").Value = 0.54999999999999993

# Row 129 — 12.8. (continued)
$ws.Range("C129").Value = 0.5541666666666667
$ws.Range("D129").Value = 0.57152777777777775

# Row 130 — 12.8. (continued)
$ws.Range("C130").Value = 0.65277777777777779
$ws.Range("D130").Value = 0.77916666666666667

# Row 131 — 12.8. (continued)
$ws.Range("C131").Value = 0.78819444444444453
$ws.Range("D131").Value = 0.83819444444444446

# Row 132 — 12.8. (continued)
$ws.Range("C132").Value = 0.89583333333333337
$ws.Range("D132").Value = 0.95138888888888884

# Restore view state: scroll so row 113 is at the top, select K126
[void]$ws.Activate()
$excel.ActiveWindow.ScrollRow = 113
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("K126").Select()
